# Project Plan - Pharmacy Error Tracker
#
# Commit: "Updated Project Plan. - Added 'Edit Error Submission Form' to
# Construction Phase ... - Reorganised order of Use Case completion in
# Construction Phase."
#
# The second table in the document ("Project milestones and objectives" /
# iteration schedule) has a Construction Phase block with rows C-1, C-2,
# C-3 (and C-4+). The 5th column of rows C-1/C-2/C-3 holds two paragraphs:
# the use-case summary (1st paragraph) and the testing note (2nd
# paragraph, left untouched). We only rewrite the 1st paragraph's text in
# each of those three cells.

function Set-CellFirstParagraphText {
    param(
        $Doc,
        $Table,
        [int]$Row,
        [int]$Col,
        [string]$NewText
    )

    $cell = $Table.Cell($Row, $Col)
    $cellRange = $cell.Range
    $cellStart = $cellRange.Start
    $cellText = $cellRange.Text

    # The first paragraph ends right before its paragraph mark (Cr, chr 13).
    $crIndex = $cellText.IndexOf([char]13)
    if ($crIndex -lt 0) {
        $paraRange = $cellRange
    } else {
        $paraRange = $Doc.Range($cellStart, $cellStart + $crIndex)
    }

    $paraRange.Text = $NewText
}

$d = $word.ActiveDocument
$scheduleTable = $d.Tables.Item(2)

# C-1 (row 2): add "Edit Error Submission Form", drop "Send Report to
# Contacts" / "Add a Contact" in favour of adding "Output Error Data to
# Excel".
Set-CellFirstParagraphText $d $scheduleTable 2 5 `
    'Implement supporting use cases "Edit Error Submission Form", "Modify Error in System", and "Output Error Data to Excel"'

# C-2 (row 3): drop "Output Error Data to Excel" / "Edit a Contact",
# gaining "Manage User Details – Preferences", "Change Password", and
# "Delete a User".
Set-CellFirstParagraphText $d $scheduleTable 3 5 `
    'Implement use cases "Add a User", "Manage User Details – Preferences", "Change Password", and "Delete a User"'

# C-3 (row 4): now groups the contact-related use cases together.
Set-CellFirstParagraphText $d $scheduleTable 4 5 `
    'Implement use cases "Add a Contact", "Edit a Contact", "Send Report to Contacts", and "Remove a Contact"'
